$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 68.42105263157895
$ws.Range("D2").Value = 74.3421052631579
$ws.Range("C3").Value = 55.26315789473685
$ws.Range("D3").Value = 79.60526315789474
$ws.Range("C4").Value = 55.26315789473685
$ws.Range("D4").Value = 67.76315789473685
$ws.Range("C5").Value = 60.52631578947368
$ws.Range("C18").Value = 0.631578947368421
$ws.Range("D18").Value = 71.71052631578947
$ws.Range("C19").Value = 0.631578947368421
$ws.Range("D19").Value = 71.71052631578947
$ws.Range("C20").Value = 0.6578947368421053
$ws.Range("D20").Value = 71.71052631578947
$ws.Range("C21").Value = 0.6052631578947368
$ws.Range("D21").Value = 78.94736842105263
$ws.Range("C22").Value = 0.6052631578947368
$ws.Range("D22").Value = 78.94736842105263
$ws.Range("C23").Value = 0.6052631578947368
$ws.Range("D23").Value = 78.94736842105263
$ws.Range("C24").Value = 0.5526315789473685
$ws.Range("D24").Value = 78.94736842105263
$ws.Range("C25").Value = 0.5526315789473685
$ws.Range("D25").Value = 78.94736842105263
$ws.Range("C26").Value = 0.5526315789473685
$ws.Range("D26").Value = 78.94736842105263
$ws.Range("C27").Value = 0.7105263157894737
$ws.Range("D27").Value = 75.6578947368421
$ws.Range("C28").Value = 0.7105263157894737
$ws.Range("D28").Value = 75.6578947368421
$ws.Range("C29").Value = 0.7105263157894737
$ws.Range("D29").Value = 75.6578947368421
$ws.Range("C72").Value = 0.6842105263157895
$ws.Range("C73").Value = 0.6578947368421053
$ws.Range("C74").Value = 0.5789473684210527
$ws.Range("C75").Value = 0.6052631578947368
$ws.Range("C84").Value = 0.6842105263157895
$ws.Range("D84").Value = 0.7763157894736842
$ws.Range("D85").Value = 0.8026315789473685
$ws.Range("C86").Value = 0.631578947368421
$ws.Range("D86").Value = 0.7763157894736842
$ws.Range("C87").Value = 0.7105263157894737
$ws.Range("D87").Value = 0.75
$ws.Range("C108").Value = 0.6578947368421053
$ws.Range("D108").Value = 0.9539473684210527
$ws.Range("C109").Value = 0.5789473684210527
$ws.Range("C110").Value = 0.4473684210526316
$ws.Range("D110").Value = 0.9736842105263158
$ws.Range("C111").Value = 0.7105263157894737
$ws.Range("D111").Value = 0.9342105263157895
$ws.Range("C124").Value = 60.52631578947368
$ws.Range("D124").Value = 69.73684210526315
$ws.Range("C125").Value = 76.31578947368422
$ws.Range("D125").Value = 73.02631578947368
$ws.Range("D126").Value = 64.47368421052632
$ws.Range("D127").Value = 70.39473684210526
$ws.Range("C140").Value = 63.1578947368421
$ws.Range("D140").Value = 0.8618421052631579
$ws.Range("C141").Value = 55.26315789473685
$ws.Range("D141").Value = 0.8947368421052632
$ws.Range("C142").Value = 71.05263157894737
$ws.Range("D142").Value = 0.8486842105263158
$ws.Range("C143").Value = 73.68421052631578
$ws.Range("D143").Value = 0.8157894736842105
$ws.Range("C156").Value = 63.1578947368421
$ws.Range("D156").Value = 0.8289473684210527
$ws.Range("C157").Value = 55.26315789473685
$ws.Range("D157").Value = 0.8618421052631579
$ws.Range("C158").Value = 63.1578947368421
$ws.Range("D158").Value = 0.7894736842105263
$ws.Range("D159").Value = 0.7894736842105263
